$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.424.94"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +6.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.811.25"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +6.11%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.48"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3834"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.07%  "
$ws.Range("E8").Value = "  +3.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3520"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.233"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07745"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.55"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +12.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.615"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.222"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +5.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.813.39"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +6.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001122"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06772"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "86.86"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +6.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.0000"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.78"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +9.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.536"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +7.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.13"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.399.14"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.472"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.700"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +8.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.32"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +16.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.491"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +14.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "154.57"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.015.15"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.71%  "
$ws.Range("E31").Value = "  +7.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.361"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +7.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.093"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.80"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +7.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08815"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.725"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.629"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.7054"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +15.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06544"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2263"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +6.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02410"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.016"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.260"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.93"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6586"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +12.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.050"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.182"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +8.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.60"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.94%  "
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.54"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.35%  "
